# Started work on live chat plugin.
#
# The task list ("Tasks" sheet) had a not-yet-started entry in row 151:
# "Apply AdminLTE Theme (...) as admin theme for main KoreCMS solution."
# / "Will be better than just a plain, ugly theme". That task is removed
# by deleting the whole row, which shifts every row below it up by one
# (row numbers printed in column A shift down by one, and the final
# trailing blank row disappears because the sheet is one row shorter).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")
$ws.Activate()

# Deleting the entire row removes the task and pulls rows 152:229 up to
# 151:228, which also renumbers the shifted A-column counters and trims
# the sheet's used range from H229 down to H228.
$ws.Rows("151:151").Delete()

# Keep the selection on the cell shown in the updated sheet.
$ws.Range("E152").Select()
